# kwb-variables: restrict the variable list to the ones actually used for
# the outcome plot / correlation / linear-model / Moran's I steps.
# A handful of rows that were previously flagged Include=1 turned out not
# to be needed, so they are flipped to 0 before re-applying the "Include"
# AutoFilter (column D) to show only the rows where Include = 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$noLongerIncluded = @(69, 70, 71, 72, 73, 74, 77, 80, 81, 82, 83)
foreach ($r in $noLongerIncluded) {
    $ws.Range("D$r").Value = 0
}

# Re-apply the AutoFilter on the "Include" column (D, the 4th column of
# the A1:F118 range) so only Include = 1 rows remain visible.
$rng = $ws.Range("A1:F118")
[void]$rng.AutoFilter(4, @("1"), 7)

# Matches the author's final on-screen selection after filtering.
[void]$ws.Range("A66:A68").Select()
